$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.8655728697776794
$ws.Range("B1").Value = 1.345322847366333
$ws.Range("C1").Value = 4.710809230804443
$ws.Range("D1").Value = 3.526508808135986
$ws.Range("E1").Value = 0.4903959929943085
